{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The commit removes an entire trailing block from the \"Requisitos\" section:\n//   - the blank paragraph right after \"LOQ4037: Qu\u00edmica Org\u00e2nica I (Requisito fraco)\"\n//   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n// The \"LOQ4037...\" paragraph itself, and the blank/page-break paragraphs that\n// follow the removed block, are left untouched.\nconst items = paragraphs.items;\n\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"LOQ4037: Qu\u00edmica Org\u00e2nica I (Requisito fraco)\") {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  const isBlank = (p) => !!p && p.text.trim() === \"\";\n  const isJupiter = (p) =>\n    !!p && p.text.trim() === \"Ver no Jupiter Salvar em pdf Salvar em docx\";\n  const isCopyright = (p) =>\n    !!p && p.text.trim().startsWith(\"\u00a9 2020 . Contact: luizeleno@usp.br.\");\n\n  const toDelete = [];\n  const next1 = items[anchorIndex + 1];\n  const next2 = items[anchorIndex + 2];\n  const next3 = items[anchorIndex + 3];\n\n  if (isBlank(next1)) toDelete.push(next1);\n  if (isJupiter(next2)) toDelete.push(next2);\n  if (isCopyright(next3)) toDelete.push(next3);\n\n  // Delete highest index first so earlier live references stay valid even if\n  // the host doesn't auto-adjust queued sibling proxies.\n  toDelete\n    .sort((a, b) => items.indexOf(b) - items.indexOf(a))\n    .forEach((p) => p.delete());\n\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# The commit removes an entire trailing block from the \"Requisitos\" section:\n#   - the blank paragraph right after \"LOQ4037: Qu\u00edmica Org\u00e2nica I (Requisito fraco)\"\n#   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n# The \"LOQ4037...\" paragraph itself, and the blank/page-break paragraphs that\n# follow the removed block, are left untouched.\n\n# Locate the \"LOQ4037: ...\" paragraph that anchors the block to remove.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $ptext = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\")\n    if ($ptext.StartsWith(\"LOQ4037: Qu\")) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ne -1) {\n    $count = $d.Paragraphs.Count\n\n    $idxBlank = $anchorIndex + 1\n    $idxJupiter = $anchorIndex + 2\n    $idxCopyright = $anchorIndex + 3\n\n    $blankOk = ($idxBlank -le $count) -and ($d.Paragraphs.Item($idxBlank).Range.Text.TrimEnd(\"`r\").Trim() -eq \"\")\n    $jupiterOk = ($idxJupiter -le $count) -and ($d.Paragraphs.Item($idxJupiter).Range.Text.TrimEnd(\"`r\").Trim() -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\")\n    $copyrightOk = ($idxCopyright -le $count) -and ($d.Paragraphs.Item($idxCopyright).Range.Text.TrimEnd(\"`r\").StartsWith([char]0xA9 + \" 2020 . Contact: luizeleno@usp.br.\"))\n\n    # Delete from the highest index down so earlier (not-yet-deleted) indices\n    # keep pointing at the correct paragraphs.\n    if ($copyrightOk) {\n        $d.Paragraphs.Item($idxCopyright).Range.Delete()\n    }\n    if ($jupiterOk) {\n        $d.Paragraphs.Item($idxJupiter).Range.Delete()\n    }\n    if ($blankOk) {\n        $d.Paragraphs.Item($idxBlank).Range.Delete()\n    }\n}\n"}
